# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6967
$ws1.Range("F7").Value  = 6840
$ws1.Range("F10").Value = 1287
$ws1.Range("F22").Value = 165
$ws1.Range("F23").Value = 633
$ws1.Range("F25").Value = 233

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6967
$ws4.Range("F7").Value  = 6840
$ws4.Range("F17").Value = 48
$ws4.Range("F22").Value = 46
$ws4.Range("F24").Value = 165
$ws4.Range("F25").Value = 633
$ws4.Range("F27").Value = 233
